# Update the "想去人数" (want-to-go count) figures in column F for the rows
# that changed between crawls, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 2269
    3 = 1713
    5 = 1089
    6 = 818
    8 = 5833
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
